$d = $word.ActiveDocument

# 1. Change highlight color from red to green for the "(10 Puntos) En su implementación..." paragraph
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*cuatro métodos recursivos*") {
        $p.Range.HighlightColorIndex = 11  # wdGreen
    }
}

Write-Host "done"
